# Refresh the crypto price/rank table to the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. All source values in this sheet are
# stored as plain text (even the numeric-looking ones), so each value is written
# with a leading apostrophe to force text entry, then the style is reset back to
# "Normal" so no stray number-format/quote-prefix style is left on the cell.
$updates = @{
    'D2' = '244.26'
    'G2' = '10'
    'D3' = '25.02'
    'G3' = '10'
    'B4' = 'HuobiToken'
    'C4' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D4' = '5.185'
    'E4' = '3HuobiTokenHT'
    'G4' = '10'
    'B5' = 'Cronos'
    'C5' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D5' = '0.05751'
    'E5' = '4CronosCRO'
    'G5' = '10'
    'B6' = 'KuCoinToken'
    'C6' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'D6' = '6.473'
    'E6' = '5KuCoinTokenKCS'
    'G6' = '10'
    'B7' = 'GateToken'
    'C7' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D7' = '3.112'
    'E7' = '6GateTokenGT'
    'G7' = '10'
    'B8' = 'MXToken'
    'C8' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D8' = '0.8119'
    'E8' = '7MXTokenMX'
    'G8' = '10'
    'B9' = 'FTXToken'
    'C9' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D9' = '0.8386'
    'E9' = '8FTXTokenFTT'
    'G9' = '10'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D10' = '0.1339'
    'E10' = '9WazirXWRX'
    'G10' = '10'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D11' = '0.06957'
    'E11' = '10MandalaExchangeTokenMDX'
    'G11' = '10'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D12' = '0.02844'
    'E12' = '11BitrueCoinBTR'
    'G12' = '10'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D13' = '0.09370'
    'E13' = '12BitMartTokenBMX'
    'G13' = '10'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D14' = '0.001509'
    'E14' = '13BitForexTokenBF'
    'G14' = '10'
    'B15' = 'One'
    'C15' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'D15' = '0.0005998'
    'E15' = '14OneONE'
    'G15' = '10'
    'D16' = '0.006074'
    'G16' = '10'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.501'
    'E17' = '16LEOLEO'
    'G17' = '10'
    'B18' = 'BTSEToken'
    'C18' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D18' = '2.109'
    'E18' = '17BTSETokenBTSE'
    'G18' = '10'
    'D19' = '0.3173'
    'G19' = '10'
    'G20' = '10'
    'D21' = '0.1300'
    'G21' = '10'
    'D22' = '3.743'
    'G22' = '10'
    'D23' = '0.04660'
    'G23' = '10'
    'G24' = '10'
    'D25' = '0.001236'
    'G25' = '10'
    'D26' = '0.004269'
    'G26' = '10'
    'D27' = '0.00008697'
    'G27' = '10'
    'D28' = '0.0002313'
    'G28' = '10'
    'G29' = '10'
    'G30' = '10'
    'G31' = '10'
    'G32' = '10'
    'G33' = '10'
    'G34' = '10'
    'G35' = '10'
    'G36' = '10'
    'G37' = '10'
    'G38' = '10'
    'G39' = '10'
    'D40' = '0.03612'
    'G40' = '10'
    'D41' = '0.006358'
    'E41' = '40KickTokenKICKBestin24h'
    'G41' = '10'
    'D42' = '0.1050'
    'G42' = '10'
    'D43' = '0.002999'
    'G43' = '10'
    'D44' = '0.007332'
    'G44' = '10'
    'D45' = '0.00005298'
    'G45' = '10'
    'G46' = '10'
    'D47' = '0.2699'
    'G47' = '10'
    'D48' = '0.002274'
    'G48' = '10'
    'G49' = '10'
    'G50' = '10'
    'G51' = '10'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
